$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-03 Sunday" "2024-03-04 Monday"

Replace-Text "94×46=" "40×34="
Replace-Text "78×30=" "38×49="
Replace-Text "34×41=" "80×86="
Replace-Text "64×89=" "52×27="
Replace-Text "12×88=" "67×75="
Replace-Text "73×55=" "98×63="
Replace-Text "29×61=" "56×30="
Replace-Text "41×87=" "21×90="
Replace-Text "65×60=" "23×88="
Replace-Text "16×93=" "58×77="
Replace-Text "76×17=" "68×82="
Replace-Text "74×92=" "58×47="
Replace-Text "68×99=" "24×52="
Replace-Text "89×51=" "23×55="
Replace-Text "94×35=" "63×23="
Replace-Text "48×11=" "34×75="
Replace-Text "99×62=" "89×34="
Replace-Text "18×23=" "67×78="
Replace-Text "80×98=" "37×92="
Replace-Text "66×48=" "53×65="
Replace-Text "40×81=" "62×33="
Replace-Text "30×89=" "28×45="
Replace-Text "96×55=" "14×52="
Replace-Text "98×21=" "33×59="
Replace-Text "80×14=" "69×33="
